$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.508429169654846
$ws.Range("B1").Value = 3.989296197891235
$ws.Range("C1").Value = 3.553694248199463
$ws.Range("D1").Value = 1.481488466262817
$ws.Range("E1").Value = 0.9731921553611755
